# Update "Presentación Trabajo Fin De Grado ... .pptx"
#  1) Slide 1  : fix "Francisco Leon" -> "Francisco León" (merge the two
#                runs into a single run and drop the stray endParaRPr).
#  2) Slide 2  : index entry "Tecnología Utilizada" -> "Tecnologías Utilizadas"
#  3) Slide 10 : title "Tecnología Utilizada" -> "Tecnologías Utilizadas"
#  4) Slide 11 : title "Tecnología Utilizada" -> "Tecnologías Utilizadas"

$p = $ppt.ActivePresentation

# --- Slide 1: "Francisco Leon" -> "Francisco León" -------------------------
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange

# Clearing the whole text range first and re-typing it as a single string
# (paragraphs separated by carriage returns) makes the engine rebuild each
# paragraph from scratch as one run carrying the original run's formatting,
# instead of diffing character-by-character against the old "Leon" run
# (which would keep it split in two runs plus a leftover endParaRPr).
$tr1.Text = ""
$tr1 = $sh1.TextFrame.TextRange
$tr1.Text = "David Corredor Miguel" + [char]13 + "Antonio Castillo" + [char]13 + "Francisco León"

# --- Slide 2: index entry "Tecnología Utilizada" ---------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(7)
$span2 = $tr2.Characters($para2.Start, $para2.Length)
$span2.Text = "Tecnologías Utilizadas"

# --- Slide 10: title "Tecnología Utilizada" ---------------------------------
$s10  = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(1)
$sh10.TextFrame.TextRange.Text = "Tecnologías Utilizadas"

# --- Slide 11: title "Tecnología Utilizada" ---------------------------------
$s11  = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(1)
$sh11.TextFrame.TextRange.Text = "Tecnologías Utilizadas"
